# merged evaluate.py and wordseq_test.py
# Rename the "seq0"/"seq1" row labels (col A, rows 3-12) to "seqv0"/"seqv1"
# so they read "seqv0+w2v_..." / "seqv1+w2v_..." instead of "seq0+w2v_..." /
# "seq1+w2v_...". These labels are shared strings, so the charts that plot
# Sheet1!$A$3:$A$12 as their category axis pick the new text up automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value2
    if ($v -ne $null) {
        $s = $v.ToString()
        if ($s.StartsWith("seq0")) {
            $cell.Value2 = "seqv0" + $s.Substring(4)
        } elseif ($s.StartsWith("seq1")) {
            $cell.Value2 = "seqv1" + $s.Substring(4)
        }
    }
}

# Matches the author's updated selection/scroll position in the saved view.
$ws.Range("A13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
